# "Comment and disable special arrangements"
#
# The faculty-availability matrix had one extra, non-standard interview
# slot ("Tuesday, Feb. 11 3:30 - 4:00") inserted between the regular
# 3:15-3:45 and 3:45-4:15 slots. That special arrangement is removed:
# the whole sheet row for that time slot is deleted so every later row
# (including the trailing "W" summary row) shifts up by one, the sheet's
# used range shrinks from A1:AW15 to A1:AW14, and the now-unused shared
# string for that slot disappears automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 13 is the "Tuesday, Feb. 11 3:30 - 4:00" row - delete it outright
# (not just clear its contents) so everything below shifts up.
$ws.Rows.Item(13).Delete()

# Leave the selection on the row that now occupies row 13 (matches the
# author re-selecting the whole row after removing it).
$ws.Rows.Item(13).Select()
